$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "10 GB" column (G) figure for the Xeon W-10885M row (row 3): 54.46 -> 54.67
$ws.Range("G3").Value = 54.67

# Move/restore the active selection to G4 (was G8)
$ws.Range("G4").Select()
